$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header changes ---
$ws.Range("B1").Value = "Timestamp"
$ws.Range("C1").Value = "Nome Completo"

# --- Per-row data: Timestamp + Name lists ---
$timestamps = @(
    "11/7/2017 15:03:28",
    "11/7/2017 15:06:29",
    "11/7/2017 15:17:47",
    "11/7/2017 15:20:46",
    "11/7/2017 15:21:31",
    "11/7/2017 15:21:45",
    "11/7/2017 15:22:11",
    "11/7/2017 15:22:45",
    "11/7/2017 15:24:19",
    "11/7/2017 15:28:36",
    "11/7/2017 15:29:33",
    "11/7/2017 15:31:34",
    "11/7/2017 15:32:36",
    "11/7/2017 15:37:08",
    "11/7/2017 15:55:06",
    "11/7/2017 16:01:21",
    "11/7/2017 16:07:06",
    "11/7/2017 17:12:54",
    "11/7/2017 17:39:56",
    "11/7/2017 18:04:14",
    "11/7/2017 18:29:10",
    "11/7/2017 18:37:57",
    "11/7/2017 23:01:51",
    "11/8/2017 8:54:17",
    "11/8/2017 9:50:05"
)

$names = @(
    "David Fogelman",
    "Wesley Gabriel Albano da Silva",
    "André Neustein",
    "Luca Salimon Nascimento",
    "Manoela Cirne lima de campos",
    "Pedro Villas Boas Dias",
    "Rodrigo coutinho",
    "Emanuelle Moço",
    "Giovana Lemes ",
    "gabriel pizzighini salvador ",
    "Luca Ribeiro Noto",
    "Pedro Cunial",
    "Luigi crespi corradi",
    "Iago Rainha Mendes",
    "Guilherme Benavente Chicarelle",
    "Bruna D'Urso de Oliveira",
    "Natália De Munno Farah",
    "Lucas Sozzi de Jesus",
    "Juliana Costa Pessoa",
    "Ana Capriles ",
    "Barbara Freire",
    "Lucas Chen Alba",
    "Bruno Arthur Cesconetto",
    "João Gabriel Rodrigues Edivirges ",
    "Matteo Iannoni"
)

for ($i = 0; $i -lt $timestamps.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $timestamps[$i]
    $ws.Cells.Item($row, 3).Value = $names[$i]
}

# --- Remove hyperlinks that used to live on column C ---
$ws.Hyperlinks.Delete()

# --- Clear the old hyperlink styling (blue/underline) from column C ---
$ws.Range("C2:C26").Style = "Normal"

